$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 6
$ws.Range("C23").Value = 4
$ws.Range("C24").Select()
